$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook and name it.
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "ODI Batting Extra"

# Header row: reuse the existing bold/bordered/centered header style from the
# "ODI Batting" sheet (rather than building a fresh style) so the style sheet
# stays as close as possible to the original workbook's.
$headerSrc = $wb.Worksheets.Item("ODI Batting")
$headerSrc.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats

$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# MATCH_CODE, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL and MAN_OF_MATCH are stored as TEXT
# (even when they look numeric); BATTING_POSITION is stored as a real number when present.
$data = @(
    @("4276", "", "", "", "", "NO"),
    @("4277", 2, "0", "0", "", "NO"),
    @("4300", 2, "1", "0", "1.68%", "NO"),
    @("4376", 2, "10", "0", "24.75%", "YES"),
    @("4432", "", "", "", "", "NO"),
    @("4433", "", "", "", "", "NO")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $ws.Cells.Item($r + 2, $c + 1)
        $value = $row[$c]
        if ($c -eq 1) {
            # BATTING_POSITION: a real number when present, otherwise blank text.
            if ($value -ne "") {
                $cell.Value = $value
            } else {
                $cell.NumberFormat = "@"
                $cell.Value = "'"
                $cell.Style = "Normal"
            }
        } else {
            # Force text storage, then drop back to the default ("Normal") style
            # so the written cell carries no explicit style index, matching the
            # un-styled data rows in the original sheets. A lone leading
            # apostrophe produces an empty *text* cell (as opposed to leaving
            # the cell truly blank) for entries with no value.
            $cell.NumberFormat = "@"
            if ($value -eq "") {
                $cell.Value = "'"
            } else {
                $cell.Value = $value
            }
            $cell.Style = "Normal"
        }
    }
}
